# Add two new authority rows to the end of the sheet, mirroring the
# formatting of the existing data rows (row 105 in particular).
# Shared strings get created column-first (both names, then both viaf
# urls) to match how the data was typed in (column A down, then column B
# down), so write it in that order too.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$lastRow = 105
$r1 = $lastRow + 1
$r2 = $lastRow + 2

# Column A first (both rows), then column B (both rows) -- matches the
# order new shared-string entries were appended in the target file.
$ws.Cells.Item($r1, 1).Value = "Oberthur, Franz"
$ws.Cells.Item($r2, 1).Value = "Mai, Angelo"
$ws.Cells.Item($r1, 2).Value = "http://viaf.org/viaf/79150290"
$ws.Cells.Item($r2, 2).Value = "http://viaf.org/viaf/100180720"

# Match formatting of the row directly above by copying its formats only,
# so the existing style entries are reused instead of new ones minted.
$ws.Cells.Item($lastRow, 1).Copy()
$ws.Cells.Item($r1, 1).PasteSpecial(-4122)
$ws.Cells.Item($r2, 1).PasteSpecial(-4122)

$ws.Cells.Item($lastRow, 2).Copy()
$ws.Cells.Item($r1, 2).PasteSpecial(-4122)
$ws.Cells.Item($r2, 2).PasteSpecial(-4122)

$excel.CutCopyMode = 0

$ws.Rows.Item($r1).RowHeight = $ws.Rows.Item($lastRow).RowHeight
$ws.Rows.Item($r2).RowHeight = $ws.Rows.Item($lastRow).RowHeight

# Update selection / view to match the new bottom of the table.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 100
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B$r2").Select()
